$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the bold/bordered header style (already used by B1/C1/D1) to A1 as well,
# since A1 becomes a header cell ("VideoName") in the new layout.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Re-write the header row, shifted one column to the left (A:C instead of B:D)
$ws.Range("A1").Value = "VideoName"
$ws.Range("B1").Value = "Views"
$ws.Range("C1").Value = "Date"

# Re-write the data row: drop the old numeric index that used to sit in A2,
# move everything one column left, and use the new video's data.
$ws.Range("A2").Value = "A Single Math Equation Makes This Possible"
$ws.Range("B2").Value = 188378
$ws.Range("C2").Value = "11 Mar 2022"

# The old column D is no longer used - clear its values and formatting entirely.
$ws.Range("D1:D2").Clear() | Out-Null
